$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ------------------------------------------------------------------
# 1. Move the "ASTM A288 E:" / spring-rate reference block from row 21
#    down to row 26 (it is being pushed down to make room for two more
#    pairs of spring rows - 17/18 and 20/21).
# ------------------------------------------------------------------
$ws.Range("A21:C21").Cut($ws.Range("A26"))
$excel.CutCopyMode = $false

# Row 21 (A21:C21) is about to be completely repurposed - strip any
# left-over formatting from the cut so the new cells start from a clean
# slate (the cut leaves ghost style references behind on the source).
$ws.Range("A21:C21").ClearFormats()

# ------------------------------------------------------------------
# 2. Update all existing formulas that referenced $B$21 so they now
#    reference $B$26.
# ------------------------------------------------------------------
$ws.Range("E2").Formula  = '=(($B$26*1000000000)*(C2/1000)^4)/(64*(F2/1000)*A2)'
$ws.Range("E3").Formula  = '=(($B$26*1000000000)*(C3/1000)^4)/(64*(F3/1000)*A3)'
$ws.Range("F5").Formula  = '=((($B$26*1000000000)*(C5/1000)^4)/(64*A5*E5))*1000'
$ws.Range("F6").Formula  = '=((($B$26*1000000000)*(C6/1000)^4)/(64*A6*E6))*1000'
$ws.Range("A8").Formula  = '=(($B$26*1000000000)*(C8/1000)^4)/(64*(F8/1000)*E8)'
$ws.Range("A9").Formula  = '=(($B$26*1000000000)*(C9/1000)^4)/(64*(F9/1000)*E9)'
$ws.Range("F11").Formula = '=((($B$26*1000000000)*(C11/1000)^4)/(64*A11*E11))*1000'
$ws.Range("F12").Formula = '=((($B$26*1000000000)*(C12/1000)^4)/(64*A12*E12))*1000'
$ws.Range("A14").Formula = '=(($B$26*1000000000)*(C14/1000)^4)/(64*(F14/1000)*E14)'
$ws.Range("A15").Formula = '=(($B$26*1000000000)*(C15/1000)^4)/(64*(F15/1000)*E15)'

# ------------------------------------------------------------------
# 3. New rows 17 & 18 - "third iteration" spring calculations, pair 1.
# ------------------------------------------------------------------
$ws.Range("A17").Value = 40
$ws.Range("B17").Value = "Music Wire"
$ws.Range("C17").Formula = '=CONVERT(D17,"in","mm")'
$ws.Range("D17").Formula = '=8/32'
$ws.Range("E17").Value = 4.125
$ws.Range("F17").Formula = '=((($B$26*1000000000)*(C17/1000)^4)/(64*A17*E17))*1000'
$ws.Range("G17").Formula = '=F17-C17'

$ws.Range("A18").Value = 19
$ws.Range("B18").Value = "Music Wire"
$ws.Range("C18").Formula = '=CONVERT(D18,"in","mm")'
$ws.Range("D18").Formula = '=7/32'
$ws.Range("E18").Value = 5.125
$ws.Range("F18").Formula = '=((($B$26*1000000000)*(C18/1000)^4)/(64*A18*E18))*1000'
$ws.Range("G18").Formula = '=F18-C18'

# Formatting: columns C/D/F/G on rows 17-18 mirror the style family used
# by the earlier "9/32" (row 5/8) and "7/32" (row 12) entries.
$ws.Range("D5").Copy()
$ws.Range("D17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D12").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C6").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F5").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F6").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G5").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G6").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. New rows 20 & 21 - "third iteration" spring calculations, pair 2.
# ------------------------------------------------------------------
$ws.Range("A20").Formula = '=(($B$26*1000000000)*(C20/1000)^4)/(64*(F20/1000)*E20)'
$ws.Range("B20").Value = "Music Wire"
$ws.Range("C20").Formula = '=CONVERT(D20,"in","mm")'
$ws.Range("D20").Value = 0.25
$ws.Range("E20").Value = 4.125
$ws.Range("F20").Formula = '=G20+C20'
$ws.Range("G20").Value = 25.25
$ws.Range("I20").Value = "third iteration of calculations - add an eighth of a coil "

$ws.Range("A21").Formula = '=(($B$26*1000000000)*(C21/1000)^4)/(64*(F21/1000)*E21)'
$ws.Range("B21").Value = "Music Wire"
$ws.Range("C21").Formula = '=CONVERT(D21,"in","mm")'
$ws.Range("D21").Value = 0.21875
$ws.Range("E21").Value = 5.125
$ws.Range("F21").Formula = '=G21+C21'
$ws.Range("G21").Value = 25.25

# Formatting: columns A/C/D/F/G on rows 20-21 mirror the style family
# used by the earlier "notes" entries (rows 8/9, 14/15).
$ws.Range("A9").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C8").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C9").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D15").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F9").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F15").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G9").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G15").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 5. Re-assert the literal values/formulas one more time in case any
#    PasteSpecial(formats) step above disturbed them.
# ------------------------------------------------------------------
$ws.Range("D17").Formula = '=8/32'
$ws.Range("D18").Formula = '=7/32'
$ws.Range("D20").Value = 0.25
$ws.Range("D21").Value = 0.21875
$ws.Range("G20").Value = 25.25
$ws.Range("G21").Value = 25.25

# ------------------------------------------------------------------
# 6. Cursor/selection tracked UI state ends up on E21.
# ------------------------------------------------------------------
$ws.Range("E21").Select()

$wb.Save()
